# Applies the "all single filter scripts in CTDC" edit:
#  - inserts a new "TabName" column (B) identifying which tab (Cases/Files) a query applies to
#  - adds a third row for a new "FilesTab" query
#  - replaces the old combined Stat query text with new per-tab Cases/Stat/Files query bodies
#  - drops the now-unused "...Neo4jData.xlsx" shared string
#  - resizes/repositions the sheet view accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- query text bodies (single-quoted here-strings => no interpolation/escaping) ----

$caseQueryOld = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s WHERE f.file_format IN [''bam'']  RETURN DISTINCT coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

$casesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
WHERE f.file_format IN ['bam'] 
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE f.file_format IN ['bam'] 
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
 WHERE f.file_format IN ['bam'] 
 WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# Here-strings keep a trailing newline; strip it so the value exactly matches the target text.
$casesQuery = $casesQuery.TrimEnd("`r", "`n")
$statQuery  = $statQuery.TrimEnd("`r", "`n")
$filesQuery = $filesQuery.TrimEnd("`r", "`n")

# ---- write the brand-new shared strings first, in first-use order, so the rebuilt  ----
# ---- shared string table lands in the same order as the target workbook.           ----
$ws.Range("B1").Value = "TabName"
$ws.Range("B2").Value = "CasesTab"
$ws.Range("B3").Value = "FilesTab"
$ws.Range("C2").Value = $casesQuery
$ws.Range("D2").Value = $statQuery
$ws.Range("C3").Value = $filesQuery
$ws.Range("D3").Value = $statQuery

# ---- re-affirm the surviving cells in their new positions (content unchanged) ----
$ws.Range("A1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

$ws.Range("A2").Value = $caseQueryOld
$ws.Range("E2").Value = "TC02_Trials_Filter_AssocFileFormat-Bam_WebData.xlsx"
$ws.Range("E3").Value = "TC02_Trials_Filter_AssocFileFormat-Bam_WebData.xlsx"

# ---- styles: reset cells that inherited formatting from cells that used to live at ----
# ---- the same address (e.g. old B2 was a wrapped query cell) back to the default,   ----
# ---- then (re)apply wrap text to the cells that hold multi-line query text.         ----
$ws.Range("B1").Style = "Normal"
$ws.Range("B2").Style = "Normal"
$ws.Range("C1").Style = "Normal"
$ws.Range("D1").Style = "Normal"

$ws.Range("A2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Range("D2").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("D3").WrapText = $true

# ---- row heights to fit the (now much longer) wrapped text ----
$ws.Rows.Item(2).RowHeight = 188.5
$ws.Rows.Item(3).RowHeight = 409.5

# ---- column widths ----
$ws.Columns.Item(1).ColumnWidth = 75.81640625
$ws.Columns.Item(2).ColumnWidth = 8.81640625
$ws.Columns.Item(3).ColumnWidth = 75.81640625
$ws.Columns.Item(4).ColumnWidth = 70.26953125
$ws.Columns.Item(5).ColumnWidth = 28.54296875

# ---- sheet view: zoom, scrolled so column B is the leftmost visible column, and the ----
# ---- active selection on the new Cases-tab query cells.                              ----
$ws.Application.ActiveWindow.Zoom = 85
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C2:D3").Select()

Write-Output "done"
